$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Metadata sheet: bump the generation "Date" property
# ------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-12T09:15:29+00:00"

# ------------------------------------------------------------------
# 2. Elements sheet: add a new mapping column (AL) after the existing
#    "Mapping: RIM Mapping" column (AK) for the new draft mapping
#    "Spécification métier vers l'extension ROR HealthcareServiceSensitiveUnit"
# ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Elements")

# Clone the look (header/data styling) of the last mapping column (AK)
# onto the new column (AL) before filling in its own values.
$ws.Range("AK1:AK6").Copy()
$ws.Range("AL1:AL6").PasteSpecial(-4122)

# Header
$ws.Range("AL1").Value = "Mapping: Spécification métier vers l'extension ROR HealthcareServiceSensitiveUnit"

# Data rows - only the Extension.value[x] row carries a mapping value
$ws.Range("AL2").Value = ""
$ws.Range("AL3").Value = ""
$ws.Range("AL4").Value = ""
$ws.Range("AL5").Value = ""
$ws.Range("AL6").Value = "uniteSensible"

# Size the new column to fit its (long) header text
$ws.Columns.Item(38).ColumnWidth = 86.95
